# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 10:30"

# 2) Swap the country labels between rows 213 and 214
#    (row 213 was "Montserrat", row 214 was "Islas Malvinas";
#     after the edit, row 213 becomes "Islas Malvinas" and row 214 becomes "Montserrat",
#     while each row keeps its own numeric data)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# 3) Update statistic rows with new case data
# Row 6 - India
$ws.Range("B6").Value = 2706450
$ws.Range("C6").Value = 4846
$ws.Range("D6").Value = 1978747
$ws.Range("E6").Value = 675748
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 51955

# Row 7 - Rusia
$ws.Range("B7").Value = 932493
$ws.Range("C7").Value = 4748
$ws.Range("D7").Value = 742628
$ws.Range("E7").Value = 173993
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 132
$ws.Range("H7").Value = 15872

# Row 25 - Filipinas
$ws.Range("B25").Value = 169213
$ws.Range("C25").Value = 4836
$ws.Range("D25").Value = 112861
$ws.Range("E25").Value = 53665
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 2687

# Row 33 - Israel
$ws.Range("B33").Value = 95264
$ws.Range("C33").Value = 513
$ws.Range("D33").Value = 71167
$ws.Range("E33").Value = 23399
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 698

# Row 34 - Ucrania
$ws.Range("B34").Value = 94436
$ws.Range("C34").Value = 1616
$ws.Range("D34").Value = 48925
$ws.Range("E34").Value = 43395
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 27
$ws.Range("H34").Value = 2116

# Row 47 - Polonia
$ws.Range("D47").Value = 39643
$ws.Range("E47").Value = 15751

# Row 48 - Singapur
$ws.Range("B48").Value = 55938
$ws.Range("C48").Value = 100
$ws.Range("E48").Value = 3561

# Row 53 - Barein
$ws.Range("E53").Value = 3482
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 174

# Row 71 - Australia
$ws.Range("B71").Value = 23773
$ws.Range("C71").Value = 214
$ws.Range("D71").Value = 14928
$ws.Range("E71").Value = 8407

# Row 111 - Hong Kong
$ws.Range("E111").Value = 856
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 70

# Row 131 - Estonia
$ws.Range("B131").Value = 2200
$ws.Range("C131").Value = 8
$ws.Range("D131").Value = 1990
$ws.Range("E131").Value = 147

# Row 213 - data (now labeled "Islas Malvinas")
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214 - data (now labeled "Montserrat")
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
